$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.049.95'
$ws.Range("E2").Value = '  +0.28%  '

$ws.Range("D3").Value = '1.922.70'
$ws.Range("E3").Value = '  +0.80%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.70'
$ws.Range("E5").Value = '  +0.32%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  +0.20%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4585'
$ws.Range("E7").Value = '  -0.24%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3815'
$ws.Range("E8").Value = '  -0.12%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07751'
$ws.Range("E9").Value = '  +0.26%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9786'
$ws.Range("E10").Value = '  -0.51%  '

$ws.Range("E11").Value = '  +2.30%  '

$ws.Range("D12").Value = '1.918.65'
$ws.Range("E12").Value = '  -0.06%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.954'
$ws.Range("E14").Value = '  -0.57%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07023'
$ws.Range("E15").Value = '  -0.23%  '

$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.006'
$ws.Range("E16").Value = '  +0.15%  '

$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '84.67'
$ws.Range("E17").Value = '  +0.50%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009500'
$ws.Range("E18").Value = '  -0.42%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.69'
$ws.Range("E19").Value = '  -0.21%  '

$ws.Range("E20").Value = '  +0.15%  '

$ws.Range("D21").Value = '29.062.92'
$ws.Range("E21").Value = '  +0.34%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.351'
$ws.Range("E22").Value = '  +0.41%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.04'
$ws.Range("E23").Value = '  +0.87%  '

$ws.Range("D24").Value = '2.156.92'
$ws.Range("E24").Value = '  +0.49%  '

$ws.Range("E25").Value = '  -0.76%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '157.80'
$ws.Range("E26").Value = '  +0.51%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.05'
$ws.Range("E27").Value = '  -0.63%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.591'
$ws.Range("E28").Value = '  +0.00%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '117.88'
$ws.Range("E29").Value = '  +0.19%  '

$ws.Range("E30").Value = '  -0.40%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09298'
$ws.Range("E31").Value = '  +0.39%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8576'
$ws.Range("E32").Value = '  -0.55%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.084'
$ws.Range("E33").Value = '  -0.51%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.241'
$ws.Range("E34").Value = '  -1.28%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.015'
$ws.Range("E35").Value = '  -0.01%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05682'
$ws.Range("E36").Value = '  -0.73%  '

$ws.Range("E37").Value = '  +0.15%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.004'
$ws.Range("E38").Value = '  +0.14%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02040'
$ws.Range("E39").Value = '  -0.14%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.109'
$ws.Range("E40").Value = '  +13.03%  '

$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.419'
$ws.Range("E41").Value = '  -1.13%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5494'
$ws.Range("E42").Value = '  -0.72%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1754'
$ws.Range("E43").Value = '  -0.15%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.349'
$ws.Range("E44").Value = '  +0.51%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.000002858'
$ws.Range("E45").Value = '  +10.43%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.163'
$ws.Range("E46").Value = '  +3.82%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5195'
$ws.Range("E47").Value = '  -0.43%  '

$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06926'
$ws.Range("E48").Value = '  +1.49%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '11.18'
$ws.Range("E49").Value = '  -0.87%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '110.45'
$ws.Range("E50").Value = '  -1.29%  '

$ws.Range("E51").Value = '  -1.30%  '
